$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 212 ("...Blackbeaut/Primera..."),
# pushing the existing row 212 down to row 214 and the existing row 213
# (Blackbeaut/Segunda) down to row 215. This mirrors the diff, where those
# two original rows reappear unmodified at r="214" and r="215".
$ws.Rows.Item(212).Insert()
$ws.Rows.Item(212).Insert()

# New row 212: Angeleno / Primera, $/bandeja 18 kilos granel
$ws.Cells.Item(212, 1).Value = 10
$ws.Cells.Item(212, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(212, 3).Value = "La Araucanía"
$ws.Cells.Item(212, 4).Value = 44656
$ws.Cells.Item(212, 5).Value = 9
$ws.Cells.Item(212, 6).Value = "Fruta"
$ws.Cells.Item(212, 7).Value = 100103
$ws.Cells.Item(212, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(212, 9).Value = 100103002
$ws.Cells.Item(212, 10).Value = "Ciruela"
$ws.Cells.Item(212, 11).Value = "Angeleno"
$ws.Cells.Item(212, 12).Value = "Primera"
$ws.Cells.Item(212, 13).Value = 20
$ws.Cells.Item(212, 14).Value = 12000
$ws.Cells.Item(212, 15).Value = 12000
$ws.Cells.Item(212, 16).Value = 12000
$ws.Cells.Item(212, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(212, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(212, 19).Value = 667
$ws.Cells.Item(212, 20).Value = 18

# New row 213: Angeleno / Primera, $/bins (450 kilos)
$ws.Cells.Item(213, 1).Value = 10
$ws.Cells.Item(213, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(213, 3).Value = "La Araucanía"
$ws.Cells.Item(213, 4).Value = 44656
$ws.Cells.Item(213, 5).Value = 9
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100103
$ws.Cells.Item(213, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(213, 9).Value = 100103002
$ws.Cells.Item(213, 10).Value = "Ciruela"
$ws.Cells.Item(213, 11).Value = "Angeleno"
$ws.Cells.Item(213, 12).Value = "Primera"
$ws.Cells.Item(213, 13).Value = 2
$ws.Cells.Item(213, 14).Value = 240000
$ws.Cells.Item(213, 15).Value = 240000
$ws.Cells.Item(213, 16).Value = 240000
$ws.Cells.Item(213, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(213, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(213, 19).Value = 533
$ws.Cells.Item(213, 20).Value = 450
